# views.xlsx — constraints: switch the "eq / ineq" operator for the two
# "mean" views (rows 2 and 3) from "<" to "=", and move the active
# selection to B6 (matching the author's subsequent editing position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2 / D3 hold the comparison operator used by the view-constraint
# formulas in column I. They were "<"; the author changed both to "=".
# A leading apostrophe forces Excel to store a literal text "=" instead
# of parsing it as the start of a formula (this is also why these cells
# carry the quotePrefix-centered style).
$ws.Range("D2").Value = "'="
$ws.Range("D3").Value = "'="

# Move/restore the active cell selection to B6 on the sheet.
[void]$ws.Range("B6").Select()
